# Apply updated crypto price/volume figures (refresh run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.672.40"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "3.159.41"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.77"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.76"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +16.49%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.438"
$ws.Range("E10").Value = "  +5.51%  "
$ws.Range("E11").Value = "  +3.91%  "
$ws.Range("E12").Value = "  +3.44%  "
$ws.Range("D13").Value = "3.704.04"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.89"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("E15").Value = "  +4.51%  "
$ws.Range("D16").Value = "58.730.90"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.25"
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").Value = "3.154.49"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.38"
$ws.Range("E21").Value = "  +4.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.532"
$ws.Range("E24").Value = "  +5.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.70"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.32"
$ws.Range("E28").Value = "  +13.63%  "
$ws.Range("D29").Value = "0.0₃0869"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.33"
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.06"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.12"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.27"
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.13"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.34"
$ws.Range("E37").Value = "  +4.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.96"
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("D41").Value = "2.649.05"
$ws.Range("E41").Value = "  +5.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.29"
$ws.Range("E42").Value = "  +7.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.723"
$ws.Range("E43").Value = "  +3.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.14"
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("E45").Value = "  +7.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "3.199.76"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +14.58%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.07"
$ws.Range("E51").Value = "  +1.23%  "
